$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 4
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = 6

$ws.Range("B6").Select()
